$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new "Price" text would otherwise be auto-parsed as a number by Excel
# (losing formatting such as trailing zeros or thousands-dot groupings) are first
# switched to Text format, exactly as a user re-keying these values would need to do.
$ws.Range("D4:D11").NumberFormat = "@"
$ws.Range("D13:D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20:D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26:D37").NumberFormat = "@"
$ws.Range("D39:D41").NumberFormat = "@"
$ws.Range("D43:D45").NumberFormat = "@"
$ws.Range("D47:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.345.28'
$ws.Range("E2").Value = '  +1.69%  '
$ws.Range("D3").Value = '1.845.99'
$ws.Range("E3").Value = '  +0.69%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '244.63'
$ws.Range("E5").Value = '  +0.05%  '
$ws.Range("D6").Value = '0.6901'
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("D7").Value = '1.001'
$ws.Range("D8").Value = '0.3048'
$ws.Range("D9").Value = '0.07603'
$ws.Range("E9").Value = '  -1.16%  '
$ws.Range("D10").Value = '23.37'
$ws.Range("E10").Value = '  +0.14%  '
$ws.Range("D11").Value = '0.07719'
$ws.Range("E11").Value = '  -1.17%  '
$ws.Range("D12").Value = '1.848.08'
$ws.Range("E12").Value = '  +0.72%  '
$ws.Range("D13").Value = '5.125'
$ws.Range("E13").Value = '  +0.61%  '
$ws.Range("D14").Value = '0.6893'
$ws.Range("E14").Value = '  +1.41%  '
$ws.Range("D15").Value = '90.06'
$ws.Range("E15").Value = '  -2.82%  '
$ws.Range("D16").Value = '6.291'
$ws.Range("E16").Value = '  -2.38%  '
$ws.Range("D17").Value = '29.350.96'
$ws.Range("E17").Value = '  +1.68%  '
$ws.Range("D18").Value = '0.000008226'
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("D19").Value = '2.094.91'
$ws.Range("E19").Value = '  +0.76%  '
$ws.Range("D20").Value = '235.77'
$ws.Range("E20").Value = '  -2.68%  '
$ws.Range("D21").Value = '12.65'
$ws.Range("E21").Value = '  -0.42%  '
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").Value = '7.648'
$ws.Range("E23").Value = '  +2.76%  '
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("E25").Value = '  -0.34%  '
$ws.Range("D26").Value = '8.926'
$ws.Range("E26").Value = '  +1.59%  '
$ws.Range("D27").Value = '160.06'
$ws.Range("E27").Value = '  +0.76%  '
$ws.Range("D28").Value = '18.12'
$ws.Range("E28").Value = '  -0.71%  '
$ws.Range("D29").Value = '1.525'
$ws.Range("E29").Value = '  -1.03%  '
$ws.Range("D30").Value = '4.243'
$ws.Range("E30").Value = '  +0.54%  '
$ws.Range("D31").Value = '4.126'
$ws.Range("E31").Value = '  -0.68%  '
$ws.Range("D32").Value = '1.198'
$ws.Range("E32").Value = '  +0.05%  '
$ws.Range("D33").Value = '0.05220'
$ws.Range("E33").Value = '  +2.67%  '
$ws.Range("D34").Value = '0.7715'
$ws.Range("E34").Value = '  -0.67%  '
$ws.Range("D35").Value = '1.871'
$ws.Range("E35").Value = '  +1.08%  '
$ws.Range("D36").Value = '1.144'
$ws.Range("E36").Value = '  +0.17%  '
$ws.Range("D37").Value = '2.677'
$ws.Range("E37").Value = '  -0.70%  '
$ws.Range("D38").Value = '1.306.97'
$ws.Range("E38").Value = '  +5.53%  '
$ws.Range("D39").Value = '0.01858'
$ws.Range("E39").Value = '  +0.43%  '
$ws.Range("D40").Value = '2.705'
$ws.Range("E40").Value = '  +0.43%  '
$ws.Range("D41").Value = '0.9437'
$ws.Range("E41").Value = '  -1.02%  '
$ws.Range("E42").Value = '  -2.27%  '
$ws.Range("D43").Value = '5.762'
$ws.Range("E43").Value = '  -2.96%  '
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("D45").Value = '9.678'
$ws.Range("E45").Value = '  +0.41%  '
$ws.Range("D46").Value = '1.994.13'
$ws.Range("E46").Value = '  +0.95%  '
$ws.Range("D47").Value = '0.5217'
$ws.Range("E47").Value = '  +1.15%  '
$ws.Range("D48").Value = '1.776'
$ws.Range("E48").Value = '  +1.97%  '
$ws.Range("D49").Value = '0.00000000121'
$ws.Range("E49").Value = '  -1.10%  '
$ws.Range("D50").Value = '63.22'
$ws.Range("E50").Value = '  -1.02%  '
$ws.Range("D51").Value = '0.05938'
$ws.Range("E51").Value = '  +0.84%  '
